$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 316; existing rows 316:336 shift down to 317:337.
$ws.Rows.Item(316).Insert()

# Populate the newly inserted row 316 with this week's new price record.
$ws.Cells.Item(316, 1).Value = 2
$ws.Cells.Item(316, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(316, 3).Value = "Coquimbo"
$ws.Cells.Item(316, 4).Value = 44826
$ws.Cells.Item(316, 5).Value = 4
$ws.Cells.Item(316, 6).Value = 100112021
$ws.Cells.Item(316, 7).Value = "Ají"
$ws.Cells.Item(316, 8).Value = "Americana (o)"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 440
$ws.Cells.Item(316, 11).Value = 73000
$ws.Cells.Item(316, 12).Value = 75000
$ws.Cells.Item(316, 13).Value = 74000
$ws.Cells.Item(316, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(316, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(316, 16).Value = 2960
$ws.Cells.Item(316, 17).Value = 25
$ws.Cells.Item(316, 18).Value = "Hortaliza"
